$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new year column (O) mirroring the formatting of
# column N (the previous last year column), then fill in the 2021 values.
$ws.Range("N4:N14").Copy()
$ws.Range("O4:O14").PasteSpecial(-4122)

$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 2
$ws.Range("O6").Value = "-"
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = "-"
$ws.Range("O9").Value = "-"
$ws.Range("O10").Value = "-"
$ws.Range("O11").Value = "-"
$ws.Range("O12").Value = 1
$ws.Range("O13").Value = "-"
$ws.Range("O14").Value = "-"

# Move / record the active selection as it was left after the edit.
$ws.Range("P1").Select()
